$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("视觉中国", "视觉中国", "实达集团")
    3  = @("榕基软件", "实达集团", "平潭发展")
    4  = @("易点天下", "易点天下", "榕基软件")
    5  = @("平潭发展", "久其软件", "航天发展")
    6  = @("实达集团", "榕基软件", "视觉中国")
    7  = @("久其软件", "国风新材", "国风新材")
    8  = @("国风新材", "平潭发展", "中水渔业")
    9  = @("航天发展", "航天发展", "大鹏工业")
    10 = @("中水渔业", "中水渔业", "合富中国")
    11 = @("赛微电子", "赛微电子", "久其软件")
    12 = @("蓝色光标", "新 华 都", "易点天下")
    13 = @("合富中国", "蓝色光标", "浪潮软件")
    14 = @("浪潮软件", "欢瑞世纪", "华夏幸福")
    15 = @("特发信息", "特发信息", "凯美特气")
    16 = @("新 华 都", "贵广网络", "华胜天成")
    17 = @("欢瑞世纪", "东方财富", "深中华A")
    18 = @("北新路桥", "龙溪股份", "新华都")
    19 = @("贵广网络", "浪潮软件", "海南海药")
    20 = @("深中华A", "深中华A", "蓝色光标")
    21 = @("龙溪股份", "合富中国", "欢瑞世纪")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
